{"js": "// \"fix: add space to template\"\n//\n// Adds a leading space before the two merge-field placeholders in the\n// \"Authorization category\" table rows:\n//   {d.authCat[i].value}    ->  \" {d.authCat[i].value}\"\n//   {d.authCat[i+1].value}  ->  \" {d.authCat[i+1].value}\"\n//\n// Both placeholders live in their own paragraph/table cell, so a literal\n// text search uniquely locates each one; we then insert a single space\n// immediately before the match.\n\nasync function addLeadingSpace(searchText) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find target text: ${searchText}`);\n  }\n\n  for (const hit of results.items) {\n    hit.insertText(\" \", \"Before\");\n  }\n  await context.sync();\n}\n\nawait addLeadingSpace(\"{d.authCat[i].value}\");\nawait addLeadingSpace(\"{d.authCat[i+1].value}\");\n", "ps1": "# fix: add space to template\n#\n# Adds a leading space before the two merge-field placeholders in the\n# \"Authorization category\" table rows:\n#   {d.authCat[i].value}    ->  \" {d.authCat[i].value}\"\n#   {d.authCat[i+1].value}  ->  \" {d.authCat[i+1].value}\"\n\n$d = $word.ActiveDocument\n\nfunction Add-LeadingSpace($searchText) {\n    $range = $d.Content\n    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, $true)\n    if ($found) {\n        $range.Collapse(1)  # wdCollapseStart\n        $range.InsertBefore(\" \")\n    } else {\n        Write-Output \"Could not find target text: $searchText\"\n    }\n}\n\nAdd-LeadingSpace(\"{d.authCat[i].value}\")\nAdd-LeadingSpace(\"{d.authCat[i+1].value}\")\n"}
